$wb = $excel.ActiveWorkbook

$wsShip  = $wb.Worksheets.Item(1)   # ShipDB
$wsPilot = $wb.Worksheets.Item(2)   # PilotDB
$wsCard  = $wb.Worksheets.Item(3)   # CardDB

# ==========================================================================
# PilotDB (sheet2): old layout A=Name B=Text C=Pts D=Unique E=Ship F=Faction
# new layout: A=Id B=Name C=Text D=Pts E=Unique F=Ship G=Faction
#             H=PilotSkill I=Firepower J=Agility K=Hull L=Shields M=Qty
# ==========================================================================

# Insert new Id column at the front; shifts old A..F to B..G
$wsPilot.Columns.Item(1).Insert()
$wsPilot.Range("A1").Value = "Id"
$wsPilot.Range("A1").Font.Bold = $true

# Append new headers after existing Faction column (G)
$wsPilot.Range("H1").Value = "PilotSkill"
$wsPilot.Range("H1").Font.Bold = $true
$wsPilot.Range("I1").Value = "Firepower"
$wsPilot.Range("I1").Font.Bold = $true
$wsPilot.Range("J1").Value = "Agility"
$wsPilot.Range("J1").Font.Bold = $true
$wsPilot.Range("K1").Value = "Hull"
$wsPilot.Range("K1").Font.Bold = $true
$wsPilot.Range("L1").Value = "Shields"
$wsPilot.Range("L1").Font.Bold = $true
$wsPilot.Range("M1").Value = "Qty"
$wsPilot.Range("M1").Font.Bold = $true

# Fill in the new Id column values for the existing pilots
$wsPilot.Range("A2").Value = "P0001"
$wsPilot.Range("A3").Value = "P0002"
$wsPilot.Range("A4").Value = "P0003"

# Fill in the new stat columns (PilotSkill, Firepower, Agility, Hull, Shields, Qty)
$wsPilot.Range("H2").Value = 8
$wsPilot.Range("I2").Value = 3
$wsPilot.Range("J2").Value = 2
$wsPilot.Range("K2").Value = 2
$wsPilot.Range("L2").Value = 2
$wsPilot.Range("M2").Value = 1

$wsPilot.Range("H3").Value = 7
$wsPilot.Range("I3").Value = 3
$wsPilot.Range("J3").Value = 1
$wsPilot.Range("K3").Value = 5
$wsPilot.Range("L3").Value = 5
$wsPilot.Range("M3").Value = 1

$wsPilot.Range("H4").Value = 2
$wsPilot.Range("I4").Value = 2
$wsPilot.Range("J4").Value = 3
$wsPilot.Range("K4").Value = 2
$wsPilot.Range("L4").Value = 1
$wsPilot.Range("M4").Value = 2

# ==========================================================================
# CardDB (sheet3): old layout A=Name B=Type C=Pts D=Qty
# new layout: A=Id B=Name C=Text D=Pts E=Unique F=Limited G=Type H=Qty
# ==========================================================================

# Preserve the existing data before reshaping the sheet
$cardName = @($wsCard.Range("A2").Value2, $wsCard.Range("A3").Value2, $wsCard.Range("A4").Value2, $wsCard.Range("A5").Value2)
$cardType = @($wsCard.Range("B2").Value2, $wsCard.Range("B3").Value2, $wsCard.Range("B4").Value2, $wsCard.Range("B5").Value2)
$cardPts  = @($wsCard.Range("C2").Value2, $wsCard.Range("C3").Value2, $wsCard.Range("C4").Value2, $wsCard.Range("C5").Value2)
$cardQty  = @($wsCard.Range("D2").Value2, $wsCard.Range("D3").Value2, $wsCard.Range("D4").Value2, $wsCard.Range("D5").Value2)

$wsCard.Cells.Clear()

$cardIds      = @("U0001","U0002","U0003","U0004")
$cardTexts    = @("Woooo","Pew pew","Go Wroom","Spacecoke")
$cardUniques  = @("N","N","N","N")
$cardLimiteds = @("N","N","Y","N")

# Header row (Id, Name, Pts, Type, Qty use already-known shared strings)
$wsCard.Cells.Item(1,1).Value = "Id"
$wsCard.Cells.Item(1,1).Font.Bold = $true
$wsCard.Cells.Item(1,2).Value = "Name"
$wsCard.Cells.Item(1,2).Font.Bold = $true
$wsCard.Cells.Item(1,4).Value = "Pts"
$wsCard.Cells.Item(1,4).Font.Bold = $true
$wsCard.Cells.Item(1,5).Value = "Unique"
$wsCard.Cells.Item(1,5).Font.Bold = $true
$wsCard.Cells.Item(1,7).Value = "Type"
$wsCard.Cells.Item(1,7).Font.Bold = $true
$wsCard.Cells.Item(1,8).Value = "Qty"
$wsCard.Cells.Item(1,8).Font.Bold = $true

# Data rows: Id / Name / Pts / Unique / Type / Qty first
for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $wsCard.Cells.Item($r, 1).Value = $cardIds[$i]
    $wsCard.Cells.Item($r, 2).Value = $cardName[$i]
    $wsCard.Cells.Item($r, 4).Value = $cardPts[$i]
    $wsCard.Cells.Item($r, 5).Value = $cardUniques[$i]
    $wsCard.Cells.Item($r, 7).Value = $cardType[$i]
    $wsCard.Cells.Item($r, 8).Value = $cardQty[$i]
}

# "Limited" header and column data
$wsCard.Cells.Item(1,6).Value = "Limited"
$wsCard.Cells.Item(1,6).Font.Bold = $true
for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $wsCard.Cells.Item($r, 6).Value = $cardLimiteds[$i]
}

# "Text" header and column data
$wsCard.Cells.Item(1,3).Value = "Text"
$wsCard.Cells.Item(1,3).Font.Bold = $true
for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $wsCard.Cells.Item($r, 3).Value = $cardTexts[$i]
}

# ==========================================================================
# ShipDB (sheet1): old layout A=Ship B=Qty
# new layout: A=Id B=Ship C=Qty
# ==========================================================================

$wsShip.Columns.Item(1).Insert()
$wsShip.Range("A1").Value = "Id"
$wsShip.Range("A1").Font.Bold = $true

$wsShip.Range("A2").Value  = "S0001"
$wsShip.Range("A3").Value  = "S0002"
$wsShip.Range("A4").Value  = "S0003"
$wsShip.Range("A5").Value  = "S0004"
$wsShip.Range("A6").Value  = "S0005"
$wsShip.Range("A7").Value  = "S0006"
$wsShip.Range("A8").Value  = "S0007"
$wsShip.Range("A9").Value  = "S0008"
$wsShip.Range("A10").Value = "S0009"
$wsShip.Range("A11").Value = "S0010"

# ==========================================================================
# View state: ShipDB selection B5, PilotDB selection M5, CardDB selection H6
# and CardDB becomes the active (selected) sheet/tab.
# ==========================================================================

$wsShip.Range("B5").Select()
$wsPilot.Range("M5").Select()
$wsCard.Range("H6").Select()
